# Update row 7 (Ano 2025) of the faturamento anual sheet with refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2356722.35
$ws.Range("C7").Value = -46.95744062226842
$ws.Range("D7").Value = 2401
$ws.Range("E7").Value = 2401
$ws.Range("F7").Value = 981.5586630570596
$ws.Range("G7").Value = 4.627055898765819
